{"js": "// Replace the date and each two-digit multiplication result cell with its new value.\n// The mapping mirrors the OOXML diff: every w:t run in the table (plus the date\n// heading paragraph) is swapped for a new literal string, so a straightforward\n// search-and-replace over the document body reproduces the edit exactly.\nconst replacements = [\n  [\"2023-09-17 Sunday\", \"2023-09-18 Monday\"],\n  [\"33\u00d781=2673\", \"83\u00d731=2573\"],\n  [\"84\u00d751=4284\", \"17\u00d759=1003\"],\n  [\"23\u00d796=2208\", \"27\u00d752=1404\"],\n  [\"45\u00d772=3240\", \"47\u00d737=1739\"],\n  [\"49\u00d742=2058\", \"80\u00d778=6240\"],\n  [\"86\u00d779=6794\", \"83\u00d717=1411\"],\n  [\"62\u00d735=2170\", \"56\u00d781=4536\"],\n  [\"99\u00d723=2277\", \"59\u00d721=1239\"],\n  [\"32\u00d796=3072\", \"25\u00d725=625\"],\n  [\"64\u00d728=1792\", \"56\u00d750=2800\"],\n  [\"80\u00d731=2480\", \"69\u00d797=6693\"],\n  [\"12\u00d748=576\", \"80\u00d795=7600\"],\n  [\"84\u00d734=2856\", \"65\u00d754=3510\"],\n  [\"92\u00d751=4692\", \"80\u00d780=6400\"],\n  [\"83\u00d724=1992\", \"50\u00d780=4000\"],\n  [\"37\u00d799=3663\", \"19\u00d782=1558\"],\n  [\"14\u00d755=770\", \"51\u00d723=1173\"],\n  [\"38\u00d758=2204\", \"65\u00d718=1170\"],\n  [\"51\u00d736=1836\", \"20\u00d797=1940\"],\n  [\"22\u00d714=308\", \"32\u00d738=1216\"],\n  [\"79\u00d719=1501\", \"44\u00d780=3520\"],\n  [\"34\u00d748=1632\", \"73\u00d745=3285\"],\n  [\"43\u00d798=4214\", \"72\u00d766=4752\"],\n  [\"30\u00d737=1110\", \"86\u00d799=8514\"],\n  [\"65\u00d743=2795\", \"33\u00d770=2310\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date heading and each two-digit multiplication result cell.\n# Mirrors the OOXML diff: every w:t run changes to a new literal string, so a\n# plain Find/Replace (wdReplaceAll) over the whole document body for each\n# old/new pair reproduces the edit exactly.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@('2023-09-17 Sunday', '2023-09-18 Monday')\n    ,@('33\u00d781=2673', '83\u00d731=2573')\n    ,@('84\u00d751=4284', '17\u00d759=1003')\n    ,@('23\u00d796=2208', '27\u00d752=1404')\n    ,@('45\u00d772=3240', '47\u00d737=1739')\n    ,@('49\u00d742=2058', '80\u00d778=6240')\n    ,@('86\u00d779=6794', '83\u00d717=1411')\n    ,@('62\u00d735=2170', '56\u00d781=4536')\n    ,@('99\u00d723=2277', '59\u00d721=1239')\n    ,@('32\u00d796=3072', '25\u00d725=625')\n    ,@('64\u00d728=1792', '56\u00d750=2800')\n    ,@('80\u00d731=2480', '69\u00d797=6693')\n    ,@('12\u00d748=576', '80\u00d795=7600')\n    ,@('84\u00d734=2856', '65\u00d754=3510')\n    ,@('92\u00d751=4692', '80\u00d780=6400')\n    ,@('83\u00d724=1992', '50\u00d780=4000')\n    ,@('37\u00d799=3663', '19\u00d782=1558')\n    ,@('14\u00d755=770', '51\u00d723=1173')\n    ,@('38\u00d758=2204', '65\u00d718=1170')\n    ,@('51\u00d736=1836', '20\u00d797=1940')\n    ,@('22\u00d714=308', '32\u00d738=1216')\n    ,@('79\u00d719=1501', '44\u00d780=3520')\n    ,@('34\u00d748=1632', '73\u00d745=3285')\n    ,@('43\u00d798=4214', '72\u00d766=4752')\n    ,@('30\u00d737=1110', '86\u00d799=8514')\n    ,@('65\u00d743=2795', '33\u00d770=2310')\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute(\n        $pair[0],  # FindText\n        $true,     # MatchCase\n        $true,     # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $pair[1],  # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"Text not found: $($pair[0])\"\n    }\n}\n"}
